$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.721.67"
$ws.Range("E2").Value = "  -2.02%  "

$ws.Range("D3").Value = "1.538.26"
$ws.Range("E3").Value = "  -1.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3934"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3199"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07202"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.073"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.770"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.637"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001096"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.20%  "

$ws.Range("D17").Value = "1.523.57"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06611"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9997"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.150"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.361"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("D25").Value = "21.728.56"
$ws.Range("E25").Value = "  -1.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.386"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.873"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").Value = "1.704.30"
$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.086"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9654"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08104"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.563"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.200"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.500"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02232"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05992"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2048"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.183"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5825"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.727"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5593"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.895"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.165"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06726"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.09%  "
